$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 122: add the new FTHG/FTAG/FTR (H/I/J) values and fill in the
# closing-odds columns (U:AC) that were previously placeholder zeros.
# Columns K:T are left completely untouched.
# ---------------------------------------------------------------------------
$ws.Range("H122").Value = 1
$ws.Range("I122").Value = 1
$ws.Range("J122").Value = "D"

$ws.Range("U122").Value = 1.8
$ws.Range("V122").Value = 2
$ws.Range("W122").Value = -1
$ws.Range("X122").Value = 1.875
$ws.Range("Y122").Value = -1
$ws.Range("Z122").Value = -0.5
$ws.Range("AA122").Value = 0.3875
$ws.Range("AB122").Value = 0.4
$ws.Range("AC122").Value = -0.5

# ---------------------------------------------------------------------------
# Row 123: brand new match row (Sabah vs Araz FK).
# Copy formatting from row 122's id/date cells first so the new cells pick
# up the same bold/bordered (A) and date (E) number formats, then overwrite
# with the real values.
# ---------------------------------------------------------------------------
$ws.Range("A122").Copy($ws.Range("A123"))
$ws.Range("E122").Copy($ws.Range("E123"))

$ws.Range("A123").Value = 121
$ws.Range("B123").Value = 7011612
$ws.Range("C123").Value = "Azerbaijan Premier League"
$ws.Range("D123").Value = "Azerbaijan Premier League"
$ws.Range("E123").Value = 45353.47916666666
$ws.Range("F123").Value = "Sabah"
$ws.Range("G123").Value = "Araz FK"
$ws.Range("H123").Value = 2
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = "H"
$ws.Range("K123").Value = 1.85
$ws.Range("L123").Value = 3.3
$ws.Range("M123").Value = 3.75
$ws.Range("N123").Value = 1.85
$ws.Range("O123").Value = 3.3
$ws.Range("P123").Value = 3.75
$ws.Range("Q123").Value = -0.5
$ws.Range("R123").Value = 1.9
$ws.Range("S123").Value = 1.9
$ws.Range("T123").Value = 2.25
$ws.Range("U123").Value = 2.025
$ws.Range("V123").Value = 1.775
$ws.Range("W123").Value = 0.8500000000000001
$ws.Range("X123").Value = -1
$ws.Range("Y123").Value = -1
$ws.Range("Z123").Value = 0.8999999999999999
$ws.Range("AA123").Value = -1
$ws.Range("AB123").Value = -0.5
$ws.Range("AC123").Value = 0.3875

# ---------------------------------------------------------------------------
# Row 124: brand new match row (FK Gabala vs Sabail FC). This match has no
# result yet, so H/I/J (FTHG/FTAG/FTR) and AB/AC (PL_AhOver/PL_AhUnder) stay
# empty, matching the source diff.
# ---------------------------------------------------------------------------
$ws.Range("A122").Copy($ws.Range("A124"))
$ws.Range("E122").Copy($ws.Range("E124"))

$ws.Range("A124").Value = 122
$ws.Range("B124").Value = 7011613
$ws.Range("C124").Value = "Azerbaijan Premier League"
$ws.Range("D124").Value = "Azerbaijan Premier League"
$ws.Range("E124").Value = 45354.375
$ws.Range("F124").Value = "FK Gabala"
$ws.Range("G124").Value = "Sabail FC"
$ws.Range("K124").Value = 3.6
$ws.Range("L124").Value = 3.4
$ws.Range("M124").Value = 1.833
$ws.Range("N124").Value = 2.75
$ws.Range("O124").Value = 3.2
$ws.Range("P124").Value = 2.3
$ws.Range("Q124").Value = 0.25
$ws.Range("R124").Value = 1.75
$ws.Range("S124").Value = 2.05
$ws.Range("T124").Value = 2.5
$ws.Range("U124").Value = 2.025
$ws.Range("V124").Value = 1.775
$ws.Range("W124").Value = 0
$ws.Range("X124").Value = 0
$ws.Range("Y124").Value = 0
$ws.Range("Z124").Value = 0
$ws.Range("AA124").Value = 0
